# Natmi re-run following Dr Hou's advice: the LR-pair table for Bmp4-Bmpr1a
# now includes a third sending/target cluster ("sCs") and splits the old
# "M2" target bucket into "M1"/"M2", so every row's stats are recomputed
# and the table grows from 3x4 to 3x5 sending/target combinations (rows 2-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("ECs", "Bmp4", "Bmpr1a", "ECs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 5.122124666666667, 15.366374, 0.09590749833957242, 0.09590749833957242, 44.43799135997666, 399.94192223979, 0.05363259670523619, 0.05363259670523619),
    @("ECs", "Bmp4", "Bmpr1a", "FAPs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 34.19248466666667, 102.577454, 0.6402256641145507, 0.6402256641145507, 296.6435682601766, 2669.79211434159, 0.35802169213322, 0.35802169213322),
    @("ECs", "Bmp4", "Bmpr1a", "M1", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 2, 0.6666666666666666, 0.5566749999999999, 1.670025, 0.01042327356567948, 0.01042327356567948, 4.829542514124999, 43.465882627125, 0.005828816694990115, 0.005828816694990116),
    @("ECs", "Bmp4", "Bmpr1a", "M2", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 1.480867666666667, 4.442603, 0.02772800791168297, 0.02772800791168297, 12.84755621136167, 115.628005902255, 0.0155058268802043, 0.0155058268802043),
    @("ECs", "Bmp4", "Bmpr1a", "sCs", 3, 1, 8.675694999999999, 26.027085, 0.5592117158070719, 0.5592117158070719, 3, 1, 12.054774, 36.164322, 0.2257155560685144, 0.2257155560685144, 104.58354251793, 941.2518826613699, 0.1262227833934213, 0.1262227833934213),
    @("FAPs", "Bmp4", "Bmpr1a", "ECs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 5.122124666666667, 15.366374, 0.09590749833957242, 0.09590749833957242, 32.35159891564978, 291.164390240848, 0.03904542496885512, 0.03904542496885512),
    @("FAPs", "Bmp4", "Bmpr1a", "FAPs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 34.19248466666667, 102.577454, 0.6402256641145507, 0.6402256641145507, 215.9614655738898, 1943.653190165008, 0.2606457635127966, 0.2606457635127966),
    @("FAPs", "Bmp4", "Bmpr1a", "M1", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 2, 0.6666666666666666, 0.5566749999999999, 1.670025, 0.01042327356567948, 0.01042327356567948, 3.515987504866666, 31.6438875438, 0.004243475775977615, 0.004243475775977616),
    @("FAPs", "Bmp4", "Bmpr1a", "M2", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 1.480867666666667, 4.442603, 0.02772800791168297, 0.02772800791168297, 9.353235213295111, 84.179116919656, 0.01128850059896438, 0.01128850059896438),
    @("FAPs", "Bmp4", "Bmpr1a", "sCs", 3, 1, 6.316050666666666, 18.948152, 0.407115456505913, 0.407115456505913, 3, 1, 12.054774, 36.164322, 0.2257155560685144, 0.2257155560685144, 76.13856335921599, 685.247070232944, 0.09189229164931924, 0.09189229164931924),
    @("sCs", "Bmp4", "Bmpr1a", "ECs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 5.122124666666667, 15.366374, 0.09590749833957242, 0.09590749833957242, 2.675825243864889, 24.082427194784, 0.003229476665481111, 0.003229476665481111),
    @("sCs", "Bmp4", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 34.19248466666667, 102.577454, 0.6402256641145507, 0.6402256641145507, 17.86233634978489, 160.761027148064, 0.02155820846853409, 0.02155820846853409),
    @("sCs", "Bmp4", "Bmpr1a", "M1", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 2, 0.6666666666666666, 0.5566749999999999, 1.670025, 0.01042327356567948, 0.01042327356567948, 0.2908099889333333, 2.6172899004, 0.0003509810947117447, 0.0003509810947117447),
    @("sCs", "Bmp4", "Bmpr1a", "M2", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 1.480867666666667, 4.442603, 0.02772800791168297, 0.02772800791168297, 0.7736131670275557, 6.962518503248001, 0.0009336804325142927, 0.0009336804325142926),
    @("sCs", "Bmp4", "Bmpr1a", "sCs", 2, 0.6666666666666666, 0.5224053333333334, 1.567216, 0.03367282768701513, 0.03367282768701512, 3, 1, 12.054774, 36.164322, 0.2257155560685144, 0.2257155560685144, 6.297478229728001, 56.677304067552, 0.007600481025773888, 0.007600481025773886),
)

# Write every row (A:T), growing the sheet from 9 used rows to 16 (header + 15 data rows).
$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowValues = $newData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
